$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 347 (pushes existing rows 347..451 down to 348..452)
$ws.Rows("347:347").Insert()

# Populate the newly inserted row 347 with the new record
$ws.Range("A347").Value = 4
$ws.Range("B347").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C347").Value = "Los Lagos"
$ws.Range("D347").Value = 45093
$ws.Range("E347").Value = 10
$ws.Range("F347").Value = 100112043
$ws.Range("G347").Value = "Pepino ensalada"
$ws.Range("H347").Value = "Sin especificar"
$ws.Range("I347").Value = "Primera"
$ws.Range("J347").Value = 400
$ws.Range("K347").Value = 17000
$ws.Range("L347").Value = 17500
$ws.Range("M347").Value = 17250
$ws.Range("N347").Value = "$/caja 60 unidades"
$ws.Range("O347").Value = "Región de Arica y Parinacota"
$ws.Range("P347").Value = 288
$ws.Range("Q347").Value = 60
$ws.Range("R347").Value = "Hortaliza"
